$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: (row, col, new value, forceText)
# forceText=1 means the new value is a plain decimal number written into a text
# column (Price); without forcing the Text number format first, Excel would
# auto-coerce it into a real number (e.g. "1.00" -> 1, "0.490" -> 0.49) and lose
# the exact display text. Values like "68.998.50" have multiple dots and are
# already safe as text with no forcing needed.
$updates = @(
    @(2, 4, '68.998.50', 0),
    @(2, 5, '  -2.69%  ', 0),
    @(3, 4, '3.673.94', 0),
    @(3, 5, '  -3.86%  ', 0),
    @(4, 4, '0.998', 1),
    @(4, 5, '  -0.15%  ', 0),
    @(5, 4, '679.11', 1),
    @(5, 5, '  -4.00%  ', 0),
    @(6, 4, '161.50', 1),
    @(6, 5, '  -5.17%  ', 0),
    @(7, 4, '3.675.02', 0),
    @(7, 5, '  -3.80%  ', 0),
    @(8, 4, '0.999', 1),
    @(8, 5, '  -0.13%  ', 0),
    @(9, 4, '0.490', 1),
    @(9, 5, '  -6.17%  ', 0),
    @(10, 4, '0.147', 1),
    @(10, 5, '  -8.69%  ', 0),
    @(11, 4, '7.21', 1),
    @(11, 5, '  -2.32%  ', 0),
    @(12, 4, '0.447', 1),
    @(12, 5, '  -2.15%  ', 0),
    @(13, 4, '0.0000234', 1),
    @(13, 5, '  -7.83%  ', 0),
    @(14, 2, 'Avalanche', 0),
    @(14, 3, 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', 0),
    @(14, 4, '33.26', 1),
    @(14, 5, '  -8.68%  ', 0),
    @(15, 2, 'WrappedliquidstakedEther2.0', 0),
    @(15, 3, 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', 0),
    @(15, 4, '4.297.74', 0),
    @(15, 5, '  -3.77%  ', 0),
    @(16, 4, '3.670.44', 0),
    @(16, 5, '  -4.00%  ', 0),
    @(17, 4, '69.030.45', 0),
    @(17, 5, '  -2.66%  ', 0),
    @(18, 5, '  -1.66%  ', 0),
    @(19, 4, '16.24', 1),
    @(19, 5, '  -6.17%  ', 0),
    @(20, 4, '6.55', 1),
    @(20, 5, '  -8.86%  ', 0),
    @(21, 4, '478.69', 1),
    @(21, 5, '  -3.18%  ', 0),
    @(22, 4, '9.78', 1),
    @(22, 5, '  -7.76%  ', 0),
    @(23, 4, '0.661', 1),
    @(23, 5, '  -9.51%  ', 0),
    @(24, 4, '78.87', 1),
    @(24, 5, '  -7.81%  ', 0),
    @(25, 4, '3.817.57', 0),
    @(25, 5, '  -3.96%  ', 0),
    @(26, 2, 'InternetComputer(DFINITY)', 0),
    @(26, 3, 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', 0),
    @(26, 4, '11.55', 1),
    @(26, 5, '  -4.36%  ', 0),
    @(27, 2, 'PEPE', 0),
    @(27, 3, 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', 0),
    @(27, 4, '0.0000127', 1),
    @(27, 5, '  -11.52%  ', 0),
    @(28, 2, 'Dai', 0),
    @(28, 3, 'https://coinranking.com/coin/MoTuySvg7+dai-dai', 0),
    @(28, 4, '1.00', 1),
    @(28, 5, '  +0.04%  ', 0),
    @(29, 4, '9.38', 1),
    @(29, 5, '  -11.21%  ', 0),
    @(30, 4, '1.80', 1),
    @(30, 5, '  -13.07%  ', 0),
    @(31, 4, '2.72', 1),
    @(31, 5, '  -12.01%  ', 0),
    @(32, 4, '2.10', 1),
    @(32, 5, '  -5.90%  ', 0),
    @(33, 4, '6.69', 1),
    @(33, 5, '  -9.46%  ', 0),
    @(34, 4, '0.999', 1),
    @(34, 5, '  -0.15%  ', 0),
    @(35, 5, '  -4.87%  ', 0),
    @(36, 4, '26.70', 1),
    @(36, 5, '  -8.91%  ', 0),
    @(37, 4, '3.632.68', 0),
    @(37, 5, '  -4.16%  ', 0),
    @(38, 4, '8.47', 1),
    @(38, 5, '  -7.23%  ', 0),
    @(39, 4, '6.06', 1),
    @(39, 5, '  +1.66%  ', 0),
    @(40, 4, '0.0929', 1),
    @(40, 5, '  -8.68%  ', 0),
    @(41, 5, '  -0.03%  ', 0),
    @(42, 2, 'Stacks', 0),
    @(42, 3, 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', 0),
    @(42, 4, '2.18', 1),
    @(42, 5, '  -5.47%  ', 0),
    @(43, 2, 'FirstDigitalUSD', 0),
    @(43, 3, 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', 0),
    @(43, 4, '1.00', 1),
    @(43, 5, '  -0.13%  ', 0),
    @(44, 4, '0.953', 1),
    @(44, 5, '  -9.07%  ', 0),
    @(45, 4, '161.57', 1),
    @(45, 5, '  -1.25%  ', 0),
    @(46, 4, '48.35', 1),
    @(46, 5, '  -0.98%  ', 0),
    @(47, 4, '2.85', 1),
    @(47, 5, '  -13.52%  ', 0),
    @(48, 4, '1.31', 1),
    @(48, 5, '  -3.63%  ', 0),
    @(49, 2, 'FLOKI', 0),
    @(49, 3, 'https://coinranking.com/coin/fmHk13Rqw+floki-floki', 0),
    @(49, 4, '0.000273', 1),
    @(49, 5, '  -11.85%  ', 0),
    @(50, 2, 'Bittensor', 0),
    @(50, 3, 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', 0),
    @(50, 4, '384.58', 1),
    @(50, 5, '  -10.13%  ', 0),
    @(51, 4, '7.98', 1),
    @(51, 5, '  -8.69%  ', 0)
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u[0], $u[1])
    if ($u[3] -eq 1) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u[2]
}
